$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 ("tot_excl"), pushing the existing "w" and "f"
# rows down to rows 4 and 5.
$ws.Rows.Item(3).Insert()

# Row 3: new "tot_excl" row - carries the values the "tot" row had before
# this edit (i.e. a snapshot of the old totals).
$ws.Range("A3").Value = "tot_excl"
$ws.Range("B3").Value = 0.966114154821934
$ws.Range("C3").Value = 189.172470588797
$ws.Range("D3").Value = -0.972262885184565
$ws.Range("E3").Value = -0.968406992009873
$ws.Range("F3").Value = 0.965572319977136
$ws.Range("G3").Value = -12.417996289589
$ws.Range("H3").Value = -676.715798029307
$ws.Range("I3").Value = 3.45568157505463
$ws.Range("J3").Value = -0.443910646922249
$ws.Range("K3").Value = -0.490921519323632
$ws.Range("L3").Value = 0.965572319977136
$ws.Range("M3").Value = -0.0238422095945621
$ws.Range("N3").Value = -1.44259223513687
$ws.Range("O3").Value = 3.45568157505463

# Row 2: "tot" row - refreshed/recalculated values (now that tot_excl exists).
$ws.Range("B2").Value = 0.966114154743509
$ws.Range("C2").Value = 189.172470010823
$ws.Range("D2").Value = -0.972262885185545
$ws.Range("E2").Value = -0.968406992011181
$ws.Range("F2").Value = 0.96557232001688
$ws.Range("G2").Value = -12.4179962894858
$ws.Range("H2").Value = -676.715798022042
$ws.Range("I2").Value = 3.4556815671311
$ws.Range("J2").Value = -0.443910646922209
$ws.Range("K2").Value = -0.4909215193235
$ws.Range("L2").Value = 0.96557232001688
$ws.Range("M2").Value = -0.0238422095946534
$ws.Range("N2").Value = -1.44259223514332
$ws.Range("O2").Value = 3.4556815671311

# Row 4: "w" row (formerly row 3) - refreshed/recalculated values.
$ws.Range("B4").Value = -0.547610477356539
$ws.Range("C4").Value = -8.03788502365533
$ws.Range("D4").Value = 0.961595440589842
$ws.Range("E4").Value = 0.957350573531041
$ws.Range("F4").Value = -0.670866991702172
$ws.Range("G4").Value = 0.468883384034928
$ws.Range("H4").Value = 25.5401844642189
$ws.Range("I4").Value = -0.179980680400155
$ws.Range("J4").Value = 0.962560466016732
$ws.Range("K4").Value = 0.960273815482235
$ws.Range("L4").Value = -0.670866991702172
$ws.Range("M4").Value = 0.0369247801652082
$ws.Range("N4").Value = 2.01542000584516
$ws.Range("O4").Value = -0.179980680400155

# Row 5: "f" row (formerly row 4) keeps its original values - untouched by
# the Insert() shift above, so nothing further to do here.
